$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 15, shifting existing rows 15-91 down to 16-92
$ws.Rows(15).Insert()

# Populate the newly inserted row 15 with the new data record
$ws.Range("A15").Value = 1
$ws.Range("B15").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C15").Value = "Arica y Parinacota"
$ws.Range("D15").Value = 44859
$ws.Range("E15").Value = 15
$ws.Range("F15").Value = 100112040
$ws.Range("G15").Value = "Cilantro"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 300
$ws.Range("K15").Value = 400
$ws.Range("L15").Value = 500
$ws.Range("M15").Value = 450
$ws.Range("N15").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O15").Value = "Región de Arica y Parinacota"
$ws.Range("P15").Value = 225
$ws.Range("Q15").Value = 2
$ws.Range("R15").Value = "Hortaliza"
